$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Home/Away team names change ---
# A2: "Kungsbacka" -> "Varberg"
$ws.Range("A2").Value = "Varberg"
# B2: "Trollhattan" -> "Kungsbacka" (and it should pick up the plain/unfilled
# style used elsewhere in the sheet instead of the old header-like fill style)
$ws.Range("B2").Value = "Kungsbacka"
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 6: add player "Marten Gullberg" with his stats ---
$ws.Range("A6").Value = "Marten Gullberg"

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
$vals = @(15,2,2,5,1,3,2,6,4,5,2,1,6,1,5)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $vals[$i]
}

# Match formatting of the existing data rows (row 5) for the new row
$ws.Range("A5:P5").Copy()
$ws.Range("A6:P6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
